# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
# and the BKEXToken / CEJI / KickToken row reorder, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.47'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.75%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '29.39'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-1.23%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.143'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.43%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05783'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.96%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.618'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.42%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.180'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '5.23%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8541'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.56%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8582'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.82%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1366'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2.02%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07033'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.71%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03161'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '10.01%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09367'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.09%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001535'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '1.46%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006008'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.16%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006053'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.71%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.483'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.65%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.165'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-3.43%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3200'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.60%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03313'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '1.62%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1284'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.92%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.165'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-12.94%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04128'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.81%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.88%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001225'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '1.14%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004128'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-4.30%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001209'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '2.49%'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '3.38%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03728'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.27%'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.005865'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '10.08%'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1069'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '1.20%'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002448'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '5.96%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009186'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-5.37%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005290'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '3.65%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.00%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05798'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-42.00%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002172'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-20.30%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002099'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.00%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001999'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.00%'
